$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '23.174.82'
$ws.Range('E2').Value = '  +0.13%  '
$ws.Range('D3').Value = '1.599.86'
$ws.Range('E3').Value = '  -0.15%  '
$ws.Range('E4').Value = '  -0.29%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '0.9982'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -0.38%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '302.68'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +0.40%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.3776'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  -0.10%  '
$ws.Range('B8').Value = 'OKB'
$ws.Range('C8').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '51.88'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +3.66%  '
$ws.Range('B9').Value = 'Cardano'
$ws.Range('C9').Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.3605'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -1.36%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '1.260'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -0.31%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.9979'
$ws.Range('D11').ClearFormats()
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.08110'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  -0.63%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '22.61'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -1.80%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '6.567'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -0.98%  '
$ws.Range('B15').Value = 'Chainlink'
$ws.Range('C15').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '7.385'
$ws.Range('D15').ClearFormats()
$ws.Range('B16').Value = 'ShibaInu'
$ws.Range('C16').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '0.00001249'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  -0.56%  '
$ws.Range('D17').Value = '1.600.72'
$ws.Range('E17').Value = '  -0.09%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '93.40'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +2.12%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '0.06847'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -0.21%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '17.99'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -1.40%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '6.524'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -0.80%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.9996'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -0.20%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '12.93'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -0.37%  '
$ws.Range('D24').Value = '23.181.08'
$ws.Range('E24').Value = '  +0.16%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.395'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +2.38%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '2.986'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +6.02%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '21.14'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  +0.31%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '149.57'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -1.12%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '5.216'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -1.19%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '133.75'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +1.27%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '2.397'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -0.50%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '6.809'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -1.59%  '
$ws.Range('D33').Value = '1.776.34'
$ws.Range('E33').Value = '  -0.12%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.9801'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +3.81%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.07574'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -1.92%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '10.30'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +1.48%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.02712'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -1.51%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '6.168'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -1.49%  '
$ws.Range('B39').Value = 'Algorand'
$ws.Range('C39').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.2502'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -1.42%  '
$ws.Range('B40').Value = 'Stellar'
$ws.Range('C40').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.08791'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -1.48%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.7111'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -0.08%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '1.364'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -2.02%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '12.37'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -2.87%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '15.43'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -1.10%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.6513'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  -1.72%  '
$ws.Range('B46').Value = 'NEARProtocol'
$ws.Range('C46').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '2.297'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -0.22%  '
$ws.Range('B47').Value = 'PancakeSwap'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '4.011'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +0.77%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '132.03'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -0.28%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.07936'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -0.19%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '1.205'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -0.36%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '1.216'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +2.12%  '
